$wb = $excel.ActiveWorkbook

# Rename the second sheet ("Лист2") to "OpenAccountTest" to host the shared
# data-provider values.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "OpenAccountTest"

# Populate the new common Data Provider table.
# Write column B (currency/Dollar) before column A's second row (Ivan Ivanov)
# so the shared-string table order matches: customer, currency, Dollar, Ivan Ivanov.
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("B2").Value = "Dollar"
$ws2.Range("A2").Value = "Ivan Ivanov"

# Widen column A to fit the new header/values.
$ws2.Columns.Item(1).ColumnWidth = 13.7109375

# Make OpenAccountTest the active sheet/tab, with A3 selected.
$ws2.Activate()
$ws2.Range("A3").Select()
